# Updated cryptos list on Thu Jul 18 17:40:47 UTC 2024 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns for every existing row, and
# shift the coin list in rows 34-51 up by one (USDe dropped, Cosmos appended
# at the end) by rewriting Coin/Link/Price/Volume for those rows.
#
# Note: several Price values look numeric (e.g. "0.119", "569.86") but must
# stay as literal text (matching the source sheet's inlineStr cells), so we
# temporarily force a text NumberFormat before writing those, then restore
# the cell to the default "Normal" style so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.772.34"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "3.420.15"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.421.90"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -6.19%  "
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.119"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.427"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").Value = "4.009.64"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("E16").Value = "  -7.25%  "
$ws.Range("D17").Value = "63.899.12"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "3.412.46"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.519"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000115"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("E27").Value = "  -3.29%  "
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -5.66%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.833"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.68%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.810.76"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0724"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.24%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "327.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.59%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.85%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.58%  "
